$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109-166 down to 110-167
$ws.Rows.Item(109).Insert()

# Populate the new row 109 with the new price record
$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44488
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = "Tropicales y subtropicales"
$ws.Range("I109").Value = 100108005
$ws.Range("J109").Value = "Piña"
$ws.Range("K109").Value = "Caramelo"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 160
$ws.Range("N109").Value = 19000
$ws.Range("O109").Value = 20000
$ws.Range("P109").Value = 19500
$ws.Range("Q109").Value = "`$/caja 14 unidades"
$ws.Range("R109").Value = "Ecuador"
$ws.Range("S109").Value = 1393
$ws.Range("T109").Value = 14
